# Time Log.xlsx — add two new "Coding" time-entry rows (122 & 123) to Sheet1,
# matching a new day of work logged on 2014-12-01 (Excel serial date 41974).
# Downstream SUMIF/summary formulas on Sheet2 and the SUM on Sheet1!E152
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The shared formula used all down column E:
#   =IF(AND(NOT(ISBLANK(Bn)),NOT(ISBLANK(Cn))), (Cn-Bn) * 24 - Dn/60, "")
$eFormula = "=IF(AND(NOT(ISBLANK(RC[-3])),NOT(ISBLANK(RC[-2]))), (RC[-2]-RC[-3]) * 24 - RC[-1]/60, """")"

# Row 122: 2014-12-01, 19:21 -> 20:31, 5 min interruption, Activity = Coding
$ws.Cells.Item(122, 1).Value = 41974
$ws.Cells.Item(122, 2).Value = 0.80625000000000002
$ws.Cells.Item(122, 3).Value = 0.85486111111111107
$ws.Cells.Item(122, 4).Value = 5
$ws.Cells.Item(122, 5).FormulaR1C1 = $eFormula
$ws.Cells.Item(122, 6).Value = "Coding"

# Row 123: 2014-12-01, 21:39 -> 00:22 (next day), 5 min interruption, Activity = Coding
$ws.Cells.Item(123, 1).Value = 41974
$ws.Cells.Item(123, 2).Value = 0.90208333333333324
$ws.Cells.Item(123, 3).Value = 1.0152777777777777
$ws.Cells.Item(123, 4).Value = 5
$ws.Cells.Item(123, 5).FormulaR1C1 = $eFormula
$ws.Cells.Item(123, 6).Value = "Coding"

# Reflect the author's final selection after entering this data.
$ws.Range("D124").Select()
